# Updated cryptos list on Mon Apr  3 09:30:12 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) figures for each
# coin row, and also picks up the reordering of InternetComputer(DFINITY)
# and Aptos (rows 39/40 swapped places in the refreshed ranking).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells remain text (many values look numeric, e.g. "1.000",
# "0.9995"), so Excel doesn't silently coerce them into real numbers and
# strip significant trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.329.50'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '1.810.25'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').Value = '313.31'
$ws.Range('E5').Value = '  -1.05%  '
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('D7').Value = '0.5148'
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').Value = '0.3982'
$ws.Range('E8').Value = '  +2.83%  '
$ws.Range('D9').Value = '0.07865'
$ws.Range('E9').Value = '  -5.50%  '
$ws.Range('D10').Value = '1.116'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('D11').Value = '40.93'
$ws.Range('E11').Value = '  -2.61%  '
$ws.Range('D12').Value = '6.380'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').Value = '1.000'
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('D14').Value = '20.44'
$ws.Range('E14').Value = '  -3.61%  '
$ws.Range('D15').Value = '7.355'
$ws.Range('E15').Value = '  -1.96%  '
$ws.Range('D16').Value = '1.800.36'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '92.93'
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('D18').Value = '0.00001083'
$ws.Range('E18').Value = '  -3.49%  '
$ws.Range('D19').Value = '0.06579'
$ws.Range('E19').Value = '  -1.19%  '
$ws.Range('D20').Value = '0.9992'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('E21').Value = '  -2.49%  '
$ws.Range('D22').Value = '6.031'
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D23').Value = '28.394.18'
$ws.Range('E23').Value = '  -0.56%  '
$ws.Range('D24').Value = '11.24'
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('D25').Value = '2.227'
$ws.Range('E25').Value = '  -1.66%  '
$ws.Range('D26').Value = '160.79'
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('D27').Value = '20.61'
$ws.Range('E27').Value = '  -2.68%  '
$ws.Range('D28').Value = '2.017.14'
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('D29').Value = '2.414'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '129.13'
$ws.Range('E30').Value = '  +2.30%  '
$ws.Range('D31').Value = '0.1086'
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('D32').Value = '1.055'
$ws.Range('E32').Value = '  -4.00%  '
$ws.Range('D33').Value = '5.600'
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('D34').Value = '3.663'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('D35').Value = '0.07195'
$ws.Range('E35').Value = '  -5.31%  '
$ws.Range('D36').Value = '9.128'
$ws.Range('E36').Value = '  +4.07%  '
$ws.Range('D37').Value = '0.02340'
$ws.Range('E37').Value = '  -1.80%  '
$ws.Range('D38').Value = '0.2167'
$ws.Range('E38').Value = '  -3.08%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '5.077'
$ws.Range('E39').Value = '  -4.15%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = '11.63'
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('D41').Value = '0.6226'
$ws.Range('E41').Value = '  -2.73%  '
$ws.Range('D42').Value = '0.9993'
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('D43').Value = '1.158'
$ws.Range('E43').Value = '  -3.02%  '
$ws.Range('D44').Value = '13.30'
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('D45').Value = '0.6034'
$ws.Range('E45').Value = '  -1.70%  '
$ws.Range('D46').Value = '1.309'
$ws.Range('E46').Value = '  -6.19%  '
$ws.Range('D47').Value = '3.738'
$ws.Range('E47').Value = '  -1.81%  '
$ws.Range('D48').Value = '126.04'
$ws.Range('E48').Value = '  -1.11%  '
$ws.Range('D49').Value = '1.220'
$ws.Range('E49').Value = '  +1.10%  '
$ws.Range('D50').Value = '1.944'
$ws.Range('E50').Value = '  -2.88%  '
$ws.Range('D51').Value = '0.06858'
$ws.Range('E51').Value = '  -1.83%  '
